$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: biscuit (A) / 餅乾 (B) / bánh quy (C)
# Set in the same order the original author typed them (Vietnamese, Chinese, English)
# so the new shared-string entries land in the same table order as the target file.
$ws.Range("C21").Value = "bánh quy"
$ws.Range("B21").Value = "餅乾"
$ws.Range("A21").Value = "biscuit"

# Row 22: meat (A) / 肉類 (B) / thịt (C, mixed-font rich text)
$ws.Range("B22").Value = "肉類"
$ws.Range("A22").Value = "meat"

$c = $ws.Range("C22")
$c.Value = "thịt"
# "th" keeps the cell's default font; "ị" uses Calibri, "t" uses 微軟正黑體 Light,
# matching the multi-run rich text recorded for this shared string.
$c.Characters(3,1).Font.Name = "Calibri"
$c.Characters(3,1).Font.Size = 12
$c.Characters(4,1).Font.Name = "微軟正黑體 Light"
$c.Characters(4,1).Font.Size = 12

# The author's last recorded selection was cell C24.
$ws.Range("C24").Select()
